$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 33: update existing measurement (shared formulas in E:G recalc automatically) ---
$ws.Range("C33").Value = 516
$ws.Range("D33").Value = 548

# --- Row 34: new measurement row ---
$ws.Range("B34").Value = "5N_field--X02--Y17_0030"
$ws.Range("C34").Value = 506
$ws.Range("D34").Value = 512
$ws.Range("E34").Formula = '= (($C34-512)*0.3)^2'
$ws.Range("F34").Formula = '= (($D34-512)*0.3)^2'
$ws.Range("G34").Formula = '=SQRT($F34+$E34)'

# --- Row 35: new measurement row ---
$ws.Range("B35").Value = "6H_field--X03--Y05_0037"
$ws.Range("C35").Value = 514
$ws.Range("D35").Value = 526
$ws.Range("E35").Formula = '= (($C35-512)*0.3)^2'
$ws.Range("F35").Formula = '= (($D35-512)*0.3)^2'
$ws.Range("G35").Formula = '=SQRT($F35+$E35)'

# --- Row 37: new measurement row (string interned before row 36's, to match
#     the original author's shared-string insertion order) ---
$ws.Range("B37").Value = "7F_field--X00--Y08_0002"
$ws.Range("C37").Value = 475
$ws.Range("D37").Value = 500
$ws.Range("E37").Formula = '= (($C37-512)*0.3)^2'
$ws.Range("F37").Formula = '= (($D37-512)*0.3)^2'
$ws.Range("G37").Formula = '=SQRT($F37+$E37)'

# --- Row 36: new measurement row ---
$ws.Range("B36").Value = "6R_field--X01--Y25_0019"
$ws.Range("C36").Value = 520
$ws.Range("D36").Value = 506
$ws.Range("E36").Formula = '= (($C36-512)*0.3)^2'
$ws.Range("F36").Formula = '= (($D36-512)*0.3)^2'
$ws.Range("G36").Formula = '=SQRT($F36+$E36)'

# --- Row 41: summary stats for the new block (G32:G41) ---
$ws.Range("I41").Formula = '=AVERAGE(G32:G41)'
$ws.Range("J41").Formula = '=STDEV(G32:G41)'

# --- Restore the view: scrolled down a bit further, D38 now selected ---
$ws.Range("D38").Select()
